$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 88) with the next business day's figures
$row = 88
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "09-11-2021"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = 50000
$ws.Cells.Item($row, 3).Value = 95000
$ws.Cells.Item($row, 4).Value = 50000
$ws.Cells.Item($row, 5).Value = 40000
$ws.Cells.Item($row, 6).Value = 10000
$ws.Cells.Item($row, 7).Value = 3.23
